$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.087.44"
$ws.Range("E2").Value = "  -2.90%  "

$ws.Range("D3").Value = "2.364.74"
$ws.Range("E3").Value = "  -4.00%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.20"
$ws.Range("E5").Value = "  -2.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.66"
$ws.Range("E6").Value = "  -3.41%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -2.52%  "

$ws.Range("D9").Value = "2.368.47"
$ws.Range("E9").Value = "  -3.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.77"
$ws.Range("E12").Value = "  +2.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.324"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("D14").Value = "2.785.76"
$ws.Range("E14").Value = "  -3.86%  "

$ws.Range("D15").Value = "56.073.79"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.41"
$ws.Range("E16").Value = "  -2.98%  "

$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").Value = "2.340.60"
$ws.Range("E18").Value = "  -2.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.02"
$ws.Range("E19").Value = "  -3.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.02"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.07"
$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("E22").Value = "  -3.40%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.81"
$ws.Range("E24").Value = "  +0.69%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.369"
$ws.Range("E26").Value = "  -3.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.147"
$ws.Range("E27").Value = "  -6.24%  "

$ws.Range("E28").Value = "  -5.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.46"

$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -3.99%  "

$ws.Range("E31").Value = "  -3.54%  "

$ws.Range("E32").Value = "  +0.23%  "

$ws.Range("E33").Value = "  -6.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("E35").Value = "  -5.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("E37").Value = "  -6.10%  "

$ws.Range("E38").Value = "  -4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.12"
$ws.Range("E39").Value = "  -1.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.792"
$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.39"
$ws.Range("E41").Value = "  -5.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.36"
$ws.Range("E42").Value = "  -2.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "128.74"
$ws.Range("E43").Value = "  -6.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.73"
$ws.Range("E44").Value = "  -5.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.563"
$ws.Range("E45").Value = "  -2.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0901"
$ws.Range("E46").Value = "  -2.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "239.16"
$ws.Range("E47").Value = "  -7.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0481"
$ws.Range("E48").Value = "  -3.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0207"
$ws.Range("E49").Value = "  -4.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.01"
$ws.Range("E50").Value = "  -3.22%  "

$ws.Range("E51").Value = "  -1.09%  "
